# "Avance del flujo Emitir dictamen propuesta"
# DatosRegistrarInformeVisitaVerificacion.xlsx is a small data-source sheet
# used to drive an automated test flow. This commit advances the sample
# data: the first data row's "Cod cliente" is changed to a new test client
# code, and the (now unneeded) second data row is removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 2 ("Cod cliente" column A): swap the sample client code used by the
# test flow for a new one.
$ws.Range("A2").Value = "1940821"

# Row 3 held a second sample client (24681769) that's no longer part of the
# flow - delete the whole row so the data range shrinks to A1:H2.
$ws.Rows("3:3").Delete()

# Leave the selection where the user last clicked while editing the sheet.
$ws.Range("H6").Select()
